$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = 200
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -488
$ws.Range("N6").ClearContents()
# Row 8
$ws.Range("H8").Value = 39.75
$ws.Range("I8").Value = 39.75
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 119.25
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 19.75
# Row 15
$ws.Range("H15").Value = 533.8
$ws.Range("I15").Value = 533.8
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1601.4
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1432.4
# Row 21
$ws.Range("H21").Value = 17362.285
$ws.Range("I21").Value = 20307.2
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 20307.2
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -19839.2
$ws.Range("N21").Value = -10936
# Row 23
$ws.Range("H23").Value = 17362.285
$ws.Range("I23").Value = 20307.2
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 20307.2
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -20073.2
$ws.Range("N23").Value = -10468
# Row 28
$ws.Range("H28").Value = 1578.7727
$ws.Range("I28").Value = 1645.375
$ws.Range("J28").Value = 1401.1666
$ws.Range("K28").Value = 1645.375
$ws.Range("L28").Value = 1401.1666
$ws.Range("M28").Value = -1160.375
$ws.Range("N28").Value = -2371.1666
# Row 31
$ws.Range("H31").Value = 6860.615
$ws.Range("I31").Value = 312.57144
$ws.Range("J31").Value = 14500
$ws.Range("K31").Value = 937.71432
$ws.Range("L31").Value = 43500
$ws.Range("M31").Value = -707.71432
$ws.Range("N31").Value = -43960
# Row 76
$ws.Range("H76").Value = 3024.2424
$ws.Range("I76").Value = 3024.2424
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3024.2424
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2709.2424
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 3024.2424
$ws.Range("I79").Value = 3024.2424
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3024.2424
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1932.2424
$ws.Range("N79").ClearContents()
# Row 129
$ws.Range("H129").Value = 984.6585
$ws.Range("I129").Value = 643.1
$ws.Range("J129").Value = 1094.8387
$ws.Range("K129").Value = 1929.3
$ws.Range("L129").Value = 3284.5161
$ws.Range("M129").Value = 3070.7
$ws.Range("N129").Value = -13284.5161
# Row 132
$ws.Range("H132").Value = 1399.7632
$ws.Range("I132").Value = 1513.7742
$ws.Range("J132").Value = 894.8570999999999
$ws.Range("K132").Value = 4541.3226
$ws.Range("L132").Value = 2684.5713
$ws.Range("M132").Value = -2011.3226
$ws.Range("N132").Value = -7744.5713

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 1000000
$ws.Range("I8").Value = 1000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -999856
# Row 102
$ws.Range("H102").Value = 2050
$ws.Range("I102").Value = 2050
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2050
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -428
# Row 109
$ws.Range("H109").Value = 30100
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 30100
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 30100
$ws.Range("N109").Value = -32874
# Row 110
$ws.Range("H110").Value = 1465
$ws.Range("I110").Value = 1501.25
$ws.Range("J110").Value = 1416.6666
$ws.Range("K110").Value = 1501.25
$ws.Range("L110").Value = 1416.6666
$ws.Range("M110").Value = 543.75
$ws.Range("N110").Value = -5506.6666

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2162.8823
$ws.Range("I105").Value = 1772.3077
$ws.Range("J105").Value = 2404.6667
$ws.Range("K105").Value = 1772.3077
$ws.Range("L105").Value = 2404.6667
$ws.Range("M105").Value = -25.30770000000007
$ws.Range("N105").Value = -5898.6667
# Row 108
$ws.Range("H108").Value = 28950
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 28950
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 28950
$ws.Range("N108").Value = -36630
# Row 112
$ws.Range("H112").Value = 29887
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 29887
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 29887
$ws.Range("N112").Value = -32841

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 910
$ws.Range("I16").Value = 868.5714
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 868.5714
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -581.5714
$ws.Range("N16").Value = -1774
# Row 53
$ws.Range("H53").Value = 25666
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 25666
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 25666
$ws.Range("N53").Value = -26880
# Row 105
$ws.Range("H105").Value = 862.625
$ws.Range("I105").Value = 816.8333
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 816.8333
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 930.1667
$ws.Range("N105").Value = -4494
# Row 107
$ws.Range("H107").Value = 1336.5217
$ws.Range("I107").Value = 555
$ws.Range("J107").Value = 2052.9167
$ws.Range("K107").Value = 555
$ws.Range("L107").Value = 2052.9167
$ws.Range("M107").Value = 1365
$ws.Range("N107").Value = -5892.9167
# Row 113
$ws.Range("H113").Value = 910
$ws.Range("I113").Value = 868.5714
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 868.5714
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1301.4286
$ws.Range("N113").Value = -5540
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 138
$ws.Range("H138").Value = 35372
$ws.Range("I138").Value = 10000
$ws.Range("J138").Value = 39275.383
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 39275.383
$ws.Range("M138").Value = -4860
$ws.Range("N138").Value = -49555.383
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 11111223
$ws.Range("I7").Value = 55.375
$ws.Range("J7").Value = 20000158
$ws.Range("K7").Value = 166.125
$ws.Range("L7").Value = 60000474
$ws.Range("M7").Value = -54.125
$ws.Range("N7").Value = -60000698
# Row 80
$ws.Range("H80").Value = 1753.6666
$ws.Range("I80").Value = 1701.3334
$ws.Range("J80").Value = 1771.1111
$ws.Range("K80").Value = 5104.0002
$ws.Range("L80").Value = 5313.3333
$ws.Range("M80").Value = -4168.0002
$ws.Range("N80").Value = -7185.3333
# Row 83
$ws.Range("H83").Value = 1753.6666
$ws.Range("I83").Value = 1701.3334
$ws.Range("J83").Value = 1771.1111
$ws.Range("K83").Value = 15312.0006
$ws.Range("L83").Value = 15939.9999
$ws.Range("M83").Value = -10632.0006
$ws.Range("N83").Value = -25299.9999
# Row 100
$ws.Range("H100").Value = 3691.6667
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3691.6667
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 11075.0001
$ws.Range("N100").Value = -12697.0001
# Row 113
$ws.Range("I113").Value = 677.7778
$ws.Range("J113").Value = 708.1111
$ws.Range("K113").Value = 2033.3334
$ws.Range("L113").Value = 2124.3333
$ws.Range("M113").Value = 136.6666
$ws.Range("N113").Value = -6464.3333
# Row 131
$ws.Range("H131").Value = 891.24
$ws.Range("I131").Value = 625
$ws.Range("J131").Value = 902.3333
$ws.Range("K131").Value = 1875
$ws.Range("L131").Value = 2706.9999
$ws.Range("M131").Value = 3165
$ws.Range("N131").Value = -12786.9999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 96
$ws.Range("H96").Value = 25150
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 25150
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 25150
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -30642
# Row 113
$ws.Range("H113").Value = 1916.8
$ws.Range("I113").Value = 1263.75
$ws.Range("J113").Value = 2352.1667
$ws.Range("K113").Value = 1263.75
$ws.Range("L113").Value = 2352.1667
$ws.Range("M113").Value = 906.25
$ws.Range("N113").Value = -6692.1667
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2242.913
$ws.Range("I61").Value = 2005.5
$ws.Range("J61").Value = 2785.5715
$ws.Range("K61").Value = 2005.5
$ws.Range("L61").Value = 2785.5715
$ws.Range("M61").Value = -1803.5
$ws.Range("N61").Value = -3189.5715
# Row 100
$ws.Range("H100").Value = 5300
$ws.Range("I100").Value = 5300
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5300
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4759
$ws.Range("N100").ClearContents()
# Row 112
$ws.Range("H112").Value = 30000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 30000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
# Row 113
$ws.Range("H113").Value = 2242.913
$ws.Range("I113").Value = 2005.5
$ws.Range("J113").Value = 2785.5715
$ws.Range("K113").Value = 2005.5
$ws.Range("L113").Value = 2785.5715
$ws.Range("M113").Value = 164.5
$ws.Range("N113").Value = -7125.5715
# Row 132
$ws.Range("H132").Value = 4234.7393
$ws.Range("I132").Value = 4393.3125
$ws.Range("J132").Value = 3872.2856
$ws.Range("K132").Value = 13179.9375
$ws.Range("L132").Value = 11616.8568
$ws.Range("M132").Value = -10649.9375
$ws.Range("N132").Value = -16676.8568
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 105695
$ws.Range("I137").Value = 40390
$ws.Range("J137").Value = 171000
$ws.Range("K137").Value = 40390
$ws.Range("L137").Value = 171000
$ws.Range("M137").Value = -35290
$ws.Range("N137").Value = -181200
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
# Row 52
$ws.Range("H52").Value = 12374.25
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 12374.25
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 12374.25
$ws.Range("N52").Value = -12826.25
# Row 112
$ws.Range("H112").Value = 29825
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 29825
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 29825
$ws.Range("N112").Value = -32779
# Row 126
$ws.Range("H126").Value = 1507.3334
$ws.Range("I126").Value = 970.36365
$ws.Range("J126").Value = 1961.6923
$ws.Range("K126").Value = 2911.09095
$ws.Range("L126").Value = 5885.0769
$ws.Range("M126").Value = -441.0909499999998
$ws.Range("N126").Value = -10825.0769
